$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) contains numeric-looking text (e.g. "42.10",
# "0.560", thousands-dot prices like "63.830.05"). Plain values that parse
# as a float would otherwise get auto-converted by Excel into a real
# number (dropping trailing zeros, etc.), so those cells are temporarily
# switched to Text format while the value is written, then restored to the
# workbook default style so no stray formatting diff is introduced.

$ws.Range("D2").Value = "63.830.05"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "3.411.34"
$ws.Range("E3").Value = "  +1.54%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.99%  "
$ws.Range("D8").Value = "3.411.90"
$ws.Range("E8").Value = "  +1.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.560"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("E11").Value = "  +4.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.435"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").Value = "3.999.93"
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.134"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.21%  "
$ws.Range("E15").Value = "  +6.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.25%  "
$ws.Range("D17").Value = "63.814.65"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("D18").Value = "3.399.21"
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("E22").Value = "  -4.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  +2.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.534"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000119"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +23.09%  "
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.47%  "
$ws.Range("E31").Value = "  +6.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").Value = "2.943.72"
$ws.Range("E39").Value = "  +5.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0758"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.04%  "
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0315"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.758"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.91%  "
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +21.01%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.105"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.01%  "
